$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first")

$data = @{
    2  = @(0.3407354241124098, 0.5844350863438295, 0.4442199068852112, 0.6664982422221466, 0.5885134887136687, 19)
    3  = @(0.1726228500093137, 0.4551251828485784, 0.3113871153733175, 0.5580207123156966, 0.5460334188492619, 18)
    4  = @(0.2034452847961829, 0.3989700209739589, 0.2615574828328119, 0.5114269085928232, 0.483661004847352, 17)
    5  = @(0.3635533387143693, 0.3920234804570175, 0.2041752638615354, 0.4518575703266854, 0.2771362769876616, 16)
    6  = @(0.3372887961726861, 0.367293011109447, 0.1918266328756511, 0.4379801740668761, 0.289203773041539, 15)
    7  = @(0.3465491651461942, 0.3849987877080897, 0.2040158891459988, 0.4516811808632266, 0.3006242113637049, 14)
    8  = @(0.3728115640212805, 0.4172293592727421, 0.2174429248205509, 0.4663077576242443, 0.2915344481050866, 13)
    9  = @(0.4140746927647911, 0.4327887006433242, 0.2237006186653438, 0.4729699976376343, 0.2387302033153843, 12)
    10 = @(0.3853625486457603, 0.4009980026591189, 0.1935302347147137, 0.4399207141232539, 0.2225500727871226, 11)
    11 = @(0.3626069129865366, 0.3740635702753569, 0.1753426339092577, 0.4187393388604153, 0.2207534083381827, 10)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    $ws.Cells.Item($rowNum, 2).Value = $vals[0]
    $ws.Cells.Item($rowNum, 3).Value = $vals[1]
    $ws.Cells.Item($rowNum, 4).Value = $vals[2]
    $ws.Cells.Item($rowNum, 5).Value = $vals[3]
    $ws.Cells.Item($rowNum, 6).Value = $vals[4]
    $ws.Cells.Item($rowNum, 7).Value = $vals[5]
}
